$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# Fix E29: was stored as text "544028", needs to be numeric 544028
$ws.Range("E29").Value = 544028

# Append new rows 30 and 31 with breakout data
$ws.Range("A30").Value = "24/06/2024 04:45:33"
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = "BSE"
$ws.Range("D30").Value = "BSE (Bombay stock exchange)"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "20"
$ws.Range("F30").Value = -1.82
$ws.Range("G30").Value = 2513.5
$ws.Range("H30").Value = 289581

$ws.Range("A31").Value = "24/06/2024 04:45:33"
$ws.Range("B31").Value = 2
$ws.Range("C31").Value = "CONCOR"
$ws.Range("D31").Value = "Container Corporation Of India Limited"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "531344"
$ws.Range("F31").Value = -3.03
$ws.Range("G31").Value = 1057.9
$ws.Range("H31").Value = 869642
